$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hotel name values for Dubai and London rows
$ws.Range("B2").Value = "Jumeirah Beach Hotel"
$ws.Range("B3").Value = "Grand Plaza Apartments"

# Widen column B to fit content
$ws.Range("B2").ColumnWidth = 20.85

# Move selection to D14 as recorded in the saved view state
$ws.Range("D14").Select()
